# Auto-generated edit script: update crypto price/volume table cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''37.580.93'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '''2.069.41'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '''232.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.48%  '
$ws.Range("D6").Value = '''0.621'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.79%  '
$ws.Range("D8").Value = '''57.94'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.57%  '
$ws.Range("D9").Value = '''0.392'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.98%  '
$ws.Range("D10").Value = '''0.0803'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.69%  '
$ws.Range("E11").Value = '  -1.09%  '
$ws.Range("D12").Value = '''15.16'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.91%  '
$ws.Range("D13").Value = '''2.373.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.58%  '
$ws.Range("D14").Value = '''21.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.53%  '
$ws.Range("D15").Value = '''0.768'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.91%  '
$ws.Range("D16").Value = '''5.37'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("D17").Value = '''2.064.35'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("D18").Value = '''37.507.88'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").Value = '''6.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.96%  '
$ws.Range("D20").Value = '''70.33'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.00%  '
$ws.Range("D21").Value = '''0.0Z0838'
$ws.Range("D21").Replace("Z", [string][char]0x2083)
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.14%  '
$ws.Range("D22").Value = '''228.93'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("E24").Value = '  +0.53%  '
$ws.Range("D25").Value = '''2.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.06%  '
$ws.Range("D26").Value = '''9.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.42%  '
$ws.Range("D27").Value = '''169.58'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.00%  '
$ws.Range("D28").Value = '''0.131'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.08%  '
$ws.Range("D29").Value = '''19.36'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.97%  '
$ws.Range("D30").Value = '''1.38'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.47%  '
$ws.Range("D31").Value = '''0.121'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("D32").Value = '''4.64'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.29%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '''4.71'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.14%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '''0.0629'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.27%  '
$ws.Range("D35").Value = '''2.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.80%  '
$ws.Range("D36").Value = '''1.83'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").Value = '''3.33'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.65%  '
$ws.Range("D38").Value = '''1.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.37%  '
$ws.Range("E39").Value = '  -3.46%  '
$ws.Range("D40").Value = '''0.0225'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.65%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '''1.505.19'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.54%  '
$ws.Range("D42").Value = '''98.65'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.05%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '''1.20'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.09%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").Value = '''17.09'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.35%  '
$ws.Range("E45").Value = '  -0.85%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '''0.0952'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.25%  '
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").Value = '''4.05'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.02%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = '''1.04'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.58%  '
$ws.Range("E49").Value = '  -1.68%  '
$ws.Range("E50").Value = '  -1.76%  '
$ws.Range("D51").Value = '''2.256.60'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.77%  '
